$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New TPM-derived values for rows 2-10 (columns E..T), per updated NATMI run.
$changes = @{
    2  = @{ E=3; F=1; G=1.069299666666667; H=3.207899; I=0.003616700200628781; J=0.003616700200628781;
            M=8.908440666666667; N=26.725322; O=0.06231272032629341; P=0.06231272032629341;
            Q=9.525792635386447; R=85.73213371847801; S=0.0002253664281058305; T=0.0002253664281058305 }
    3  = @{ E=3; F=1; G=1.069299666666667; H=3.207899; I=0.003616700200628781; J=0.003616700200628781;
            O=0.3097346304939027; P=0.3097346304939027;
            Q=47.3493669772909; R=426.144302795618; S=0.001120217300248979; T=0.001120217300248979 }
    4  = @{ E=3; F=1; G=1.069299666666667; H=3.207899; I=0.003616700200628781; J=0.003616700200628781;
            M=89.774269; N=269.322807; O=0.627952649179804; P=0.627952649179804;
            Q=95.99559591694369; R=863.9603632524932; S=0.002271116472273971; T=0.002271116472273971 }
    5  = @{ I=0.8238194745364892; J=0.8238194745364891;
            M=8.908440666666667; N=26.725322; O=0.06231272032629341; P=0.06231272032629341;
            Q=2169.804807725918; R=19528.24326953326; S=0.05133443251614625; T=0.05133443251614624 }
    6  = @{ I=0.8238194745364892; J=0.8238194745364891;
            O=0.3097346304939027; P=0.3097346304939027;
            S=0.2551654205392406; T=0.2551654205392405 }
    7  = @{ I=0.8238194745364892; J=0.8238194745364891;
            M=89.774269; N=269.322807; O=0.627952649179804; P=0.627952649179804;
            Q=21866.07598063138; R=196794.6838256824; S=0.5173196214811024; T=0.5173196214811023 }
    8  = @{ G=51.01955666666666; H=153.05867; I=0.1725638252628821; J=0.1725638252628821;
            M=8.908440666666667; N=26.725322; O=0.06231272032629341; P=0.06231272032629341;
            Q=454.5046934046377; R=4090.54224064174; S=0.01075292138204134; T=0.01075292138204134 }
    9  = @{ G=51.01955666666666; H=153.05867; I=0.1725638252628821; J=0.1725638252628821;
            O=0.3097346304939027; P=0.3097346304939027;
            Q=2259.183077424215; R=20332.64769681794; S=0.05344899265441318; T=0.05344899265441318 }
    10 = @{ G=51.01955666666666; H=153.05867; I=0.1725638252628821; J=0.1725638252628821;
            M=89.774269; N=269.322807; O=0.627952649179804; P=0.627952649179804;
            Q=4580.243404454076; R=41222.19064008669; S=0.1083619112264276; T=0.1083619112264276 }
}

foreach ($row in $changes.Keys) {
    $cols = $changes[$row]
    foreach ($col in $cols.Keys) {
        $ws.Range("$col$row").Value = $cols[$col]
    }
}
